$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "TC30_Search_Typeahead" ---

# Row3: B3 "CLICK_PRE_ENTERTEXT" -> "CLICK"
$ws1.Range("B3").Value = "CLICK"

# Insert a new row 5 (shifts old row5 -> row6) and populate it with "WAIT"
$ws1.Rows("5:5").Insert()
$ws1.Range("A5:E5").Borders.LineStyle = 1
$ws1.Range("B5").Value = "WAIT"

# Update the selection / used range highlighting on sheet1
$ws1.Range("A3:XFD7").Select()

# --- Sheet2 "Testdata" ---

$ws2.Range("A5").Value = "EleType1"
$ws2.Range("B5").Value = "JSElement"
$ws2.Range("A6").Value = "EleType2"
$ws2.Range("B6").Value = "JSElement"
$ws2.Range("A5:B6").Borders.LineStyle = 1

$ws2.Activate()
$ws2.Range("A5:B6").Select()

$ws1.Activate()
